$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 848, pushing existing rows 848-887 down to 849-888.
$ws.Rows.Item(848).Insert()

# Populate the newly inserted row 848 with the new record's data.
$r = 848
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($r, 3).Value = "Bíobío"
$ws.Cells.Item($r, 4).Value = 45267
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r + 1, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = 8
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100108
$ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value = 100108006
$ws.Cells.Item($r, 10).Value = "Plátano"
$ws.Cells.Item($r, 11).Value = "Sin especificar"
$ws.Cells.Item($r, 12).Value = "Pintón"
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = 26000
$ws.Cells.Item($r, 15).Value = 26000
$ws.Cells.Item($r, 16).Value = 26000
$ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item($r, 18).Value = "Ecuador"
$ws.Cells.Item($r, 19).Value = 1300
$ws.Cells.Item($r, 20).Value = 20
